$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: "Updated 5/13" + bookmark + "/21"  ->
#           "Updated " + "6/25/21" + bookmark  (two separate runs,
#           bookmark _GoBack relocated to the very end of the line)
# -----------------------------------------------------------------

# Remove "5/13" from the first run, leaving "Updated " (8 chars) right
# before the (still collapsed) _GoBack bookmark.
$r1 = $d.Content
$r1.Find.Execute("5/13") | Out-Null
$r1.Text = ""

# Replace the trailing "/21" run (which sits right after the bookmark)
# with the new date text, in place, so it stays its own run.
$r2 = $d.Content
$r2.Find.Execute("/21") | Out-Null
$r2.Text = "6/25/21"

# The bookmark is still sitting between the two runs.  Push it past the
# new "6/25/21" run: temporarily stick an extra character after the
# paragraph's last run, recreate the bookmark collapsed right before
# that extra character (landing it exactly after "6/25/21"), then trim
# the extra character back off again.
$p1 = $d.Paragraphs(1)
$tailPos = $p1.Range.End - 1
$tmp = $d.Range($tailPos, $tailPos)
$tmp.InsertAfter("Z")

$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$newBmRange = $d.Range($tailPos, $tailPos)
$d.Bookmarks.Add("_GoBack", $newBmRange)

$bm2 = $d.Bookmarks("_GoBack")
$cleanup = $d.Range($bm2.End, $bm2.End + 1)
$cleanup.Text = ""

# -----------------------------------------------------------------
# Change 2: add two new list items after "RDMallClean_Qualtrics.csv"
# -----------------------------------------------------------------

$rFind = $d.Content
$rFind.Find.Execute("RDMallClean_Qualtrics.csv") | Out-Null
$paraCsv = $rFind.Paragraphs(1)
$endCsv = $paraCsv.Range.End
$paraCsv.Range.InsertParagraphAfter()

$newPara1Range = $d.Range($endCsv, $endCsv + 1)
$newPara1Range.InsertBefore("RDMallCleanLoss_Qualtrics.csv")

$rFind2 = $d.Content
$rFind2.Find.Execute("RDMallCleanLoss_Qualtrics.csv") | Out-Null
$paraLoss = $rFind2.Paragraphs(1)
$endLoss = $paraLoss.Range.End
$paraLoss.Range.InsertParagraphAfter()

$newPara2Range = $d.Range($endLoss, $endLoss + 1)
$newPara2Range.InsertBefore("RDMallCleanGain_Qualtrics.csv")

# -----------------------------------------------------------------
# Change 3: "code" -> "C" + "ode" + " (preprocessing only)"
#           (three separate runs)
# -----------------------------------------------------------------

$rCode = $d.Content
$rCode.Find.Execute("code", $true, $true) | Out-Null
$codeStart = $rCode.Start
$codeEnd = $rCode.End
$rCode.Text = "Code"

# Force a run split right after the "C" by adding+removing a bookmark
# at that point (this permanently splits the underlying run).
$splitPos = $codeStart + 1
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TEMP_CODE_SPLIT", $splitRange)
$bmSplit = $d.Bookmarks("TEMP_CODE_SPLIT")
$bmSplit.Delete()

# Append the suffix as its own run after "ode".
$suffixPoint = $d.Range($codeEnd, $codeEnd)
$suffixPoint.InsertAfter(" (preprocessing only)")

# -----------------------------------------------------------------
# Change 4: move <w:lastRenderedPageBreak/> from the "stateQuartiles.csv"
# run to the "RDMcleanDataColumnInfo.csv " run.
# -----------------------------------------------------------------

# Re-typing the text (even to the same value) drops the stale
# lastRenderedPageBreak marker that Word had cached on that run.
$rState = $d.Content
$rState.Find.Execute("stateQuartiles.csv") | Out-Null
$rState.Text = "stateQuartiles.csv"
